# Fgf18-Fgfr2.xlsx update: refresh NATMI LR-pair TPM results and drop the
# "Neutrophils" sending-cluster rows (cluster removed from this dataset).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Neutrophils" cluster (original data rows 14:16) is dropped entirely;
# "Resolving-Mac" (previously rows 17:19) shifts up to become rows 14:16.
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(14).Delete()

# Refresh the NATMI metric columns (E:T) for every remaining data row (2:16)
# with the newly recomputed TPM-based values.
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 1.209657333333333
$ws.Range("H2").Value2 = 3.628972
$ws.Range("I2").Value2 = 0.1420138414668746
$ws.Range("J2").Value2 = 0.1420138414668746
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 0.106124
$ws.Range("N2").Value2 = 0.318372
$ws.Range("O2").Value2 = 0.08094716512538251
$ws.Range("P2").Value2 = 0.08094716512538253
$ws.Range("Q2").Value2 = 0.1283736748426667
$ws.Range("R2").Value2 = 1.155363073584
$ws.Range("S2").Value2 = 0.011495617875309
$ws.Range("T2").Value2 = 0.011495617875309

$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 1.209657333333333
$ws.Range("H3").Value2 = 3.628972
$ws.Range("I3").Value2 = 0.1420138414668746
$ws.Range("J3").Value2 = 0.1420138414668746
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 1.092289666666667
$ws.Range("N3").Value2 = 3.276869
$ws.Range("O3").Value2 = 0.8331551016962769
$ws.Range("P3").Value2 = 0.833155101696277
$ws.Range("Q3").Value2 = 1.321296205407555
$ws.Range("R3").Value2 = 11.891665848668
$ws.Range("S3").Value2 = 0.1183195565296129
$ws.Range("T3").Value2 = 0.1183195565296129

$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 1.209657333333333
$ws.Range("H4").Value2 = 3.628972
$ws.Range("I4").Value2 = 0.1420138414668746
$ws.Range("J4").Value2 = 0.1420138414668746
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 0.1126143333333333
$ws.Range("N4").Value2 = 0.337843
$ws.Range("O4").Value2 = 0.08589773317834044
$ws.Range("P4").Value2 = 0.08589773317834046
$ws.Range("Q4").Value2 = 0.1362247541551111
$ws.Range("R4").Value2 = 1.226022787396
$ws.Range("S4").Value2 = 0.01219866706195274
$ws.Range("T4").Value2 = 0.01219866706195274

$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 5.449095666666666
$ws.Range("H5").Value2 = 16.347287
$ws.Range("I5").Value2 = 0.6397241489963273
$ws.Range("J5").Value2 = 0.6397241489963273
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 0.106124
$ws.Range("N5").Value2 = 0.318372
$ws.Range("O5").Value2 = 0.08094716512538251
$ws.Range("P5").Value2 = 0.08094716512538253
$ws.Range("Q5").Value2 = 0.5782798285293332
$ws.Range("R5").Value2 = 5.204518456763999
$ws.Range("S5").Value2 = 0.05178385632350051
$ws.Range("T5").Value2 = 0.05178385632350052

$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 5.449095666666666
$ws.Range("H6").Value2 = 16.347287
$ws.Range("I6").Value2 = 0.6397241489963273
$ws.Range("J6").Value2 = 0.6397241489963273
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 1.092289666666667
$ws.Range("N6").Value2 = 3.276869
$ws.Range("O6").Value2 = 0.8331551016962769
$ws.Range("P6").Value2 = 0.833155101696277
$ws.Range("Q6").Value2 = 5.951990889378109
$ws.Range("R6").Value2 = 53.56791800440299
$ws.Range("S6").Value2 = 0.5329894384145993
$ws.Range("T6").Value2 = 0.5329894384145994

$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 5.449095666666666
$ws.Range("H7").Value2 = 16.347287
$ws.Range("I7").Value2 = 0.6397241489963273
$ws.Range("J7").Value2 = 0.6397241489963273
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 0.1126143333333333
$ws.Range("N7").Value2 = 0.337843
$ws.Range("O7").Value2 = 0.08589773317834044
$ws.Range("P7").Value2 = 0.08589773317834046
$ws.Range("Q7").Value2 = 0.6136462757712221
$ws.Range("R7").Value2 = 5.522816481941
$ws.Range("S7").Value2 = 0.05495085425822743
$ws.Range("T7").Value2 = 0.05495085425822744

$ws.Range("E8").Value2 = 2
$ws.Range("F8").Value2 = 0.6666666666666666
$ws.Range("G8").Value2 = 0.5047176666666666
$ws.Range("H8").Value2 = 1.514153
$ws.Range("I8").Value2 = 0.05925388349609548
$ws.Range("J8").Value2 = 0.05925388349609547
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 0.106124
$ws.Range("N8").Value2 = 0.318372
$ws.Range("O8").Value2 = 0.08094716512538251
$ws.Range("P8").Value2 = 0.08094716512538253
$ws.Range("Q8").Value2 = 0.05356265765733333
$ws.Range("R8").Value2 = 0.4820639189159999
$ws.Range("S8").Value2 = 0.004796433891678618
$ws.Range("T8").Value2 = 0.004796433891678618

$ws.Range("E9").Value2 = 2
$ws.Range("F9").Value2 = 0.6666666666666666
$ws.Range("G9").Value2 = 0.5047176666666666
$ws.Range("H9").Value2 = 1.514153
$ws.Range("I9").Value2 = 0.05925388349609548
$ws.Range("J9").Value2 = 0.05925388349609547
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 1.092289666666667
$ws.Range("N9").Value2 = 3.276869
$ws.Range("O9").Value2 = 0.8331551016962769
$ws.Range("P9").Value2 = 0.833155101696277
$ws.Range("Q9").Value2 = 0.551297891884111
$ws.Range("R9").Value2 = 4.961681026957
$ws.Range("S9").Value2 = 0.04936767533008878
$ws.Range("T9").Value2 = 0.04936767533008878

$ws.Range("E10").Value2 = 2
$ws.Range("F10").Value2 = 0.6666666666666666
$ws.Range("G10").Value2 = 0.5047176666666666
$ws.Range("H10").Value2 = 1.514153
$ws.Range("I10").Value2 = 0.05925388349609548
$ws.Range("J10").Value2 = 0.05925388349609547
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 0.1126143333333333
$ws.Range("N10").Value2 = 0.337843
$ws.Range("O10").Value2 = 0.08589773317834044
$ws.Range("P10").Value2 = 0.08589773317834046
$ws.Range("Q10").Value2 = 0.05683844355322221
$ws.Range("R10").Value2 = 0.511545991979
$ws.Range("S10").Value2 = 0.005089774274328079
$ws.Range("T10").Value2 = 0.00508977427432808

$ws.Range("E11").Value2 = 3
$ws.Range("F11").Value2 = 1
$ws.Range("G11").Value2 = 1.177313666666667
$ws.Range("H11").Value2 = 3.531941
$ws.Range("I11").Value2 = 0.1382166931142909
$ws.Range("J11").Value2 = 0.1382166931142909
$ws.Range("K11").Value2 = 3
$ws.Range("L11").Value2 = 1
$ws.Range("M11").Value2 = 0.106124
$ws.Range("N11").Value2 = 0.318372
$ws.Range("O11").Value2 = 0.08094716512538251
$ws.Range("P11").Value2 = 0.08094716512538253
$ws.Range("Q11").Value2 = 0.1249412355613333
$ws.Range("R11").Value2 = 1.124471120052
$ws.Range("S11").Value2 = 0.01118824948060683
$ws.Range("T11").Value2 = 0.01118824948060683

$ws.Range("E12").Value2 = 3
$ws.Range("F12").Value2 = 1
$ws.Range("G12").Value2 = 1.177313666666667
$ws.Range("H12").Value2 = 3.531941
$ws.Range("I12").Value2 = 0.1382166931142909
$ws.Range("J12").Value2 = 0.1382166931142909
$ws.Range("K12").Value2 = 3
$ws.Range("L12").Value2 = 1
$ws.Range("M12").Value2 = 1.092289666666667
$ws.Range("N12").Value2 = 3.276869
$ws.Range("O12").Value2 = 0.8331551016962769
$ws.Range("P12").Value2 = 0.833155101696277
$ws.Range("Q12").Value2 = 1.285967552525444
$ws.Range("R12").Value2 = 11.573707972729
$ws.Range("S12").Value2 = 0.1151559430077602
$ws.Range("T12").Value2 = 0.1151559430077602

$ws.Range("E13").Value2 = 3
$ws.Range("F13").Value2 = 1
$ws.Range("G13").Value2 = 1.177313666666667
$ws.Range("H13").Value2 = 3.531941
$ws.Range("I13").Value2 = 0.1382166931142909
$ws.Range("J13").Value2 = 0.1382166931142909
$ws.Range("K13").Value2 = 3
$ws.Range("L13").Value2 = 1
$ws.Range("M13").Value2 = 0.1126143333333333
$ws.Range("N13").Value2 = 0.337843
$ws.Range("O13").Value2 = 0.08589773317834044
$ws.Range("P13").Value2 = 0.08589773317834046
$ws.Range("Q13").Value2 = 0.1325823936958889
$ws.Range("R13").Value2 = 1.193241543263
$ws.Range("S13").Value2 = 0.01187250062592393
$ws.Range("T13").Value2 = 0.01187250062592393

$ws.Range("E14").Value2 = 1
$ws.Range("F14").Value2 = 0.3333333333333333
$ws.Range("G14").Value2 = 0.177099
$ws.Range("H14").Value2 = 0.531297
$ws.Range("I14").Value2 = 0.02079143292641169
$ws.Range("J14").Value2 = 0.02079143292641169
$ws.Range("K14").Value2 = 3
$ws.Range("L14").Value2 = 1
$ws.Range("M14").Value2 = 0.106124
$ws.Range("N14").Value2 = 0.318372
$ws.Range("O14").Value2 = 0.08094716512538251
$ws.Range("P14").Value2 = 0.08094716512538253
$ws.Range("Q14").Value2 = 0.018794454276
$ws.Range("R14").Value2 = 0.169150088484
$ws.Range("S14").Value2 = 0.001683007554287562
$ws.Range("T14").Value2 = 0.001683007554287562

$ws.Range("E15").Value2 = 1
$ws.Range("F15").Value2 = 0.3333333333333333
$ws.Range("G15").Value2 = 0.177099
$ws.Range("H15").Value2 = 0.531297
$ws.Range("I15").Value2 = 0.02079143292641169
$ws.Range("J15").Value2 = 0.02079143292641169
$ws.Range("K15").Value2 = 3
$ws.Range("L15").Value2 = 1
$ws.Range("M15").Value2 = 1.092289666666667
$ws.Range("N15").Value2 = 3.276869
$ws.Range("O15").Value2 = 0.8331551016962769
$ws.Range("P15").Value2 = 0.833155101696277
$ws.Range("Q15").Value2 = 0.193443407677
$ws.Range("R15").Value2 = 1.740990669093
$ws.Range("S15").Value2 = 0.01732248841421585
$ws.Range("T15").Value2 = 0.01732248841421585

$ws.Range("E16").Value2 = 1
$ws.Range("F16").Value2 = 0.3333333333333333
$ws.Range("G16").Value2 = 0.177099
$ws.Range("H16").Value2 = 0.531297
$ws.Range("I16").Value2 = 0.02079143292641169
$ws.Range("J16").Value2 = 0.02079143292641169
$ws.Range("K16").Value2 = 3
$ws.Range("L16").Value2 = 1
$ws.Range("M16").Value2 = 0.1126143333333333
$ws.Range("N16").Value2 = 0.337843
$ws.Range("O16").Value2 = 0.08589773317834044
$ws.Range("P16").Value2 = 0.08589773317834046
$ws.Range("Q16").Value2 = 0.019943885819
$ws.Range("R16").Value2 = 0.179494972371
$ws.Range("S16").Value2 = 0.001785936957908274
$ws.Range("T16").Value2 = 0.001785936957908274

